# Update cryptocurrency price/volume snapshot (cryptos list refresh)
# Mirrors the scheduled GitHub Actions data refresh that updates
# the "Price" (column D) and "Volume(1h)" (column E) columns for
# the rows that moved since the last run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are stored as plain text (several look like numbers
# -- e.g. "22.361.91" or "1.003" -- but must stay text, matching the
# original inline-string cells). Prefixing with a leading apostrophe
# forces Excel to keep the literal text instead of coercing it to a
# number, and resetting the cell Style back to "Normal" afterwards
# drops the quote-prefix formatting so the cell keeps its original
# (unstyled) appearance.

$ws.Range("D2").Value = "'22.361.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "'1.566.67"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'290.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'0.3781"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("D8").Value = "'49.21"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.3402"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'0.07603"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").Value = "'1.138"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "'21.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'5.985"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").Value = "'6.919"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'1.567.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'0.00001133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'89.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'0.06739"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'6.195"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("D23").Value = "'11.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").Value = "'22.358.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'2.399"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "'2.689"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.22%  "
$ws.Range("D27").Value = "'20.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'147.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").Value = "'5.030"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.75%  "
$ws.Range("D30").Value = "'126.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'1.740.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "'2.015"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").Value = "'6.081"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.56%  "
$ws.Range("D34").Value = "'0.9936"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.30%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "'1.419"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.81%  "
$ws.Range("D37").Value = "'0.08459"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "'0.02509"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("D39").Value = "'0.2294"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").Value = "'0.06488"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'5.406"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.36%  "
$ws.Range("D42").Value = "'11.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.54%  "
$ws.Range("D43").Value = "'0.6319"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").Value = "'1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "'14.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").Value = "'3.806"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").Value = "'0.5935"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "'2.083"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "'1.256"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").Value = "'124.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").Value = "'0.07322"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.45%  "
